# Applies the "Updated cryptos list" price/volume refresh described in the commit diff.
# Source data cells are plain text (coinranking.com scrape), so any numeric-looking
# replacement value gets its NumberFormat forced to "@" (Text) first -- otherwise
# Excel would silently coerce strings like "2.00" or "0.0739" into numbers and drop
# the formatting (trailing/leading zeros) that the source data relies on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '35.506.71'
$ws.Range("E2").Value = '  -0.22%  '

# Row 3
$ws.Range("D3").Value = '1.922.08'
$ws.Range("E3").Value = '  +0.27%  '

# Row 4
$ws.Range("E4").Value = '  -0.55%  '

# Row 5
$ws.Range("E5").Value = '  +10.77%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '250.81'
$ws.Range("E6").Value = '  +1.44%  '

# Row 7
$ws.Range("E7").Value = '  -0.48%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '40.84'
$ws.Range("E8").Value = '  -3.27%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.358'
$ws.Range("E9").Value = '  +2.94%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.88'
$ws.Range("E10").Value = '  +7.04%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0739'
$ws.Range("E11").Value = '  +1.74%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0996'
$ws.Range("E12").Value = '  -0.58%  '

# Row 13
$ws.Range("D13").Value = '2.193.38'
$ws.Range("E13").Value = '  -0.09%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.70'
$ws.Range("E14").Value = '  +2.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.720'
$ws.Range("E15").Value = '  +2.47%  '

# Row 16
$ws.Range("D16").Value = '1.926.66'
$ws.Range("E16").Value = '  +0.47%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.94'
$ws.Range("E17").Value = '  -0.24%  '

# Row 18
$ws.Range("D18").Value = '35.529.98'
$ws.Range("E18").Value = '  -0.24%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.29'
$ws.Range("E19").Value = '  +1.09%  '

# Row 20
$ws.Range("E20").Value = '  +0.50%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.22'
$ws.Range("E21").Value = '  +4.09%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '242.45'
$ws.Range("E22").Value = '  -1.53%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.10'
$ws.Range("E23").Value = '  +4.57%  '

# Row 24
$ws.Range("E24").Value = '  -0.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("E25").Value = '  +0.57%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.39'
$ws.Range("E26").Value = '  +7.64%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.93'
$ws.Range("E27").Value = '  -2.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.74'
$ws.Range("E28").Value = '  +2.58%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.136'
$ws.Range("E29").Value = '  +5.87%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.89'
$ws.Range("E30").Value = '  +1.91%  '

# Row 31
$ws.Range("D31").Value = '4.132.80'
$ws.Range("E31").Value = '  +19.57%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.35'
$ws.Range("E32").Value = '  +3.86%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.00'
$ws.Range("E33").Value = '  +13.93%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0582'
$ws.Range("E34").Value = '  +1.47%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.30'
$ws.Range("E35").Value = '  +2.22%  '

# Row 36
$ws.Range("E36").Value = '  -0.54%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.916'
$ws.Range("E37").Value = '  -2.28%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.51'
$ws.Range("E38").Value = '  +11.58%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.10'
$ws.Range("E39").Value = '  +2.67%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.69'
$ws.Range("E40").Value = '  +11.66%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '99.99'
$ws.Range("E41").Value = '  +8.75%  '

# Row 42
$ws.Range("E42").Value = '  +2.85%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0211'
$ws.Range("E43").Value = '  -0.68%  '

# Row 44
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0657'
$ws.Range("E44").Value = '  +2.28%  '

# Row 45
$ws.Range("E45").Value = '  +2.76%  '

# Row 46
$ws.Range("D46").Value = '1.352.65'
$ws.Range("E46").Value = '  -0.52%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.42'
$ws.Range("E47").Value = '  +0.84%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("E48").Value = '  -0.45%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.68'
$ws.Range("E49").Value = '  +1.51%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.41'
$ws.Range("E50").Value = '  -3.36%  '

# Row 51
$ws.Range("E51").Value = '  -5.58%  '
